# Test data update for Oncology pops
# - NewImportLogic: swap the "Test894 - Test894" population for
#   "QOL_and_ECON - UtilityOutcome" (new name, radio button label, expected
#   source template path, and Excel/Word report filenames).
# - OldImportLogic: point the expected source template at the new
#   QOL_ECON_Testing folder (population name/labels stay "AAA - mCRPC").
# - Make "NewImportLogic" the active/selected sheet again (it was "prodfix").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("NewImportLogic")
$ws1.Range("A2").Value = "QOL_and_ECON - UtilityOutcome"
$ws1.Range("B2").Value = "QOL_and_ECON - UtilityOutcome_radio_button"
$ws1.Range("H3").Value = "ExcelReport-QOL_and_ECON - UtilityOutcome-Quality of Life-"
$ws1.Range("H4").Value = "WordReport-QOL_and_ECON - UtilityOutcome-Quality of Life-"
$ws1.Range("E2").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Testing\UtilityOutcome_Feature_Extraction_file_QOL_with_manipulated_data.xlsx"

$ws2 = $wb.Worksheets.Item("OldImportLogic")
$ws2.Range("E2").Value = "\Testdata\Templates\UtilityOutcome\QOL_ECON_Testing\OldImportLogic_QOL\QOL_OldImportExpectedResult_mainpulated_data.xlsx"
$ws2.Range("E2").Select()

$ws1.Activate()
$ws1.Range("E2").Select()
